$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.771073460578918
$ws.Range("B1").Value = 1.796207189559937
$ws.Range("C1").Value = 2.073438882827759
$ws.Range("D1").Value = 1.993209719657898
$ws.Range("E1").Value = 2.91481351852417
